$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (X_COORD_CD): add counts, taller row for wrapped description text
$ws.Range("C21").Value2 = 5384167
$ws.Range("E21").Value2 = 195869
$ws.Rows.Item(21).RowHeight = 48

# Row 22 (Y_COORD_CD): add matching counts
$ws.Range("C22").Value2 = 5384167
$ws.Range("E22").Value2 = 195869

# Widen column C to fit the new counts comfortably
$ws.Columns.Item(3).ColumnWidth = 24.1640625

# New font used for a trailing placeholder cell below the table
$r = $ws.Range("B29")
$r.Font.Name = "Menlo"
$r.Font.Size = 11
$r.Font.Color = 16777215

# Move the active selection to reflect where editing left off
$null = $ws.Range("C30").Select()
